# Auto-generated edit script: updates cryptocurrency price/volume data
# per the commit "Updated cryptos list on Sat Nov 23 23:56:30 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.647.10"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "3.392.20"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'254.87"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'650.76"
$ws.Range("E6").Value = "  +2.74%  "
$ws.Range("D7").Value = "'1.46"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'0.429"
$ws.Range("E8").Value = "  +3.95%  "
$ws.Range("D9").Value = "'1.06"
$ws.Range("E9").Value = "  +4.76%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "3.389.55"
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("E12").Value = "  +3.50%  "
$ws.Range("D13").Value = "'41.42"
$ws.Range("E13").Value = "  -4.07%  "
$ws.Range("D14").Value = "'6.32"
$ws.Range("E14").Value = "  +16.04%  "
$ws.Range("D15").Value = "'0.0000259"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "97.266.36"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "4.031.96"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").Value = "'8.47"
$ws.Range("E18").Value = "  +27.70%  "
$ws.Range("D19").Value = "3.389.61"
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("D20").Value = "'17.36"
$ws.Range("E20").Value = "  +5.01%  "
$ws.Range("D21").Value = "'0.515"
$ws.Range("E21").Value = "  +51.77%  "
$ws.Range("D22").Value = "'10.71"
$ws.Range("E22").Value = "  +8.26%  "
$ws.Range("D23").Value = "'3.43"
$ws.Range("E23").Value = "  -3.72%  "
$ws.Range("D24").Value = "'508.78"
$ws.Range("E24").Value = "  +3.75%  "
$ws.Range("D25").Value = "'0.0000205"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'99.19"
$ws.Range("E26").Value = "  +10.06%  "
$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").Value = "'6.10"
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("D28").Value = "'12.72"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("D29").Value = "3.578.57"
$ws.Range("E29").Value = "  +2.14%  "
$ws.Range("D30").Value = "'0.154"
$ws.Range("E30").Value = "  +4.97%  "
$ws.Range("D31").Value = "'0.204"
$ws.Range("E31").Value = "  +6.25%  "
$ws.Range("D32").Value = "'11.38"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").Value = "'0.998"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'0.569"
$ws.Range("E35").Value = "  +15.27%  "
$ws.Range("D36").Value = "'29.54"
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("D37").Value = "'2.26"
$ws.Range("E37").Value = "  +12.74%  "
$ws.Range("D38").Value = "'7.64"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").Value = "'1.42"
$ws.Range("E39").Value = "  +10.76%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.153"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'520.66"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "'0.854"
$ws.Range("E43").Value = "  +5.57%  "
$ws.Range("D44").Value = "'0.0421"
$ws.Range("E44").Value = "  +18.38%  "
$ws.Range("D45").Value = "'3.66"
$ws.Range("E45").Value = "  -4.61%  "
$ws.Range("D46").Value = "'3.26"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("D47").Value = "'5.37"
$ws.Range("E47").Value = "  +8.80%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.20"
$ws.Range("E49").Value = "  +8.82%  "
$ws.Range("D50").Value = "'1.59"
$ws.Range("E50").Value = "  +7.19%  "
$ws.Range("D51").Value = "'2.06"
$ws.Range("E51").Value = "  +1.00%  "

# Reset the number format/style on forced-text price cells back to the default
# "Normal" style so they match the original (unstyled) cells exactly.
$ws.Range("D2:D51").Style = "Normal"

